# "Add pdf and quiz" — appends a new closing "Quiz!" slide (with a link to a
# Google Forms quiz) right after the existing last slide ("Obrigado!").
#
# The new slide reuses the same slide layout / shape formatting as the
# "Obrigado!" slide that currently closes the deck, so the most faithful way
# to reproduce it is to duplicate that slide and then restyle its single
# textbox: reposition/resize it, and replace its text with the three
# "Quiz!" / blank / hyperlink paragraphs.

$p = $ppt.ActivePresentation

# The deck currently ends with slide 26 ("Obrigado!") - duplicate it so the
# new slide inherits its layout, textbox formatting and run properties.
$lastSlide = $p.Slides.Item($p.Slides.Count)
$newRange = $lastSlide.Duplicate()
$newSlide = $newRange.Item(1)

$shp = $newSlide.Shapes.Item(1)

# Reposition / resize the textbox (values taken from the target EMU, COM
# works in points where 1 pt = 12700 EMU).
$shp.Left = 508000 / 12700
$shp.Top = 2211163 / 12700
$shp.Width = 11175999 / 12700
$shp.Height = 2435674 / 12700

# Replace the "Obrigado!" text with the Quiz! content: a title line, a
# blank line, and a bold/italic hyperlinked line pointing at the quiz form.
$tr = $shp.TextFrame.TextRange
$tr.Text = "Quiz!" + [char]13 + [char]13 + "https://forms.gle/FRc6wmk1q6ioX7F48"

# Style the hyperlink paragraph (bold, italic) and wire the hyperlink.
$linkPara = $tr.Paragraphs(3, 1)
$linkPara.Font.Bold = -1
$linkPara.Font.Italic = -1
$linkPara.ActionSettings(1).Hyperlink.Address = "https://forms.gle/FRc6wmk1q6ioX7F48"
